$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve text formatting for the Price column so numeric-looking strings
# (e.g. "1.00", "305.38") are not auto-converted to real numbers by Excel.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "39.643.73"
$ws.Range("E2").Value = "  -4.70%  "

$ws.Range("D3").Value = "2.314.64"
$ws.Range("E3").Value = "  -5.92%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "305.38"
$ws.Range("E5").Value = "  -4.15%  "

$ws.Range("D6").Value = "83.49"
$ws.Range("E6").Value = "  -8.55%  "

$ws.Range("E7").Value = "  -3.83%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").Value = "0.478"
$ws.Range("E9").Value = "  -5.11%  "

$ws.Range("D10").Value = "0.0803"
$ws.Range("E10").Value = "  -5.86%  "

$ws.Range("D11").Value = "29.53"
$ws.Range("E11").Value = "  -9.58%  "

$ws.Range("E12").Value = "  +0.12%  "

$ws.Range("D13").Value = "2.675.47"
$ws.Range("E13").Value = "  -5.75%  "

$ws.Range("D14").Value = "6.34"
$ws.Range("E14").Value = "  -7.27%  "

$ws.Range("D15").Value = "14.52"
$ws.Range("E15").Value = "  -5.80%  "

$ws.Range("D16").Value = "2.317.75"
$ws.Range("E16").Value = "  -6.05%  "

$ws.Range("D17").Value = "0.745"
$ws.Range("E17").Value = "  -5.04%  "

$ws.Range("D18").Value = "39.670.09"

$ws.Range("E19").Value = "  -4.66%  "

$ws.Range("D20").Value = "6.03"
$ws.Range("E20").Value = "  -5.53%  "

$ws.Range("D21").Value = "67.16"
$ws.Range("E21").Value = "  -6.50%  "

$ws.Range("D22").Value = "10.47"
$ws.Range("E22").Value = "  -6.26%  "

$ws.Range("D23").Value = "234.48"
$ws.Range("E23").Value = "  -1.56%  "

$ws.Range("D24").Value = "2.52"
$ws.Range("E24").Value = "  -7.92%  "

$ws.Range("E25").Value = "  +0.19%  "

$ws.Range("E26").Value = "  -7.63%  "

$ws.Range("D27").Value = "22.87"
$ws.Range("E27").Value = "  -7.13%  "

$ws.Range("E28").Value = "  -5.43%  "

$ws.Range("D29").Value = "9.14"
$ws.Range("E29").Value = "  -5.31%  "

$ws.Range("D30").Value = "33.94"
$ws.Range("E30").Value = "  -5.94%  "

$ws.Range("D31").Value = "150.23"
$ws.Range("E31").Value = "  -4.54%  "

$ws.Range("E32").Value = "  +0.07%  "

$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  -6.74%  "

$ws.Range("E34").Value = "  -4.57%  "

$ws.Range("D35").Value = "0.0710"
$ws.Range("E35").Value = "  -6.62%  "

$ws.Range("E36").Value = "  -2.58%  "

$ws.Range("E37").Value = "  -3.96%  "

$ws.Range("E38").Value = "  -6.49%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -7.78%  "

$ws.Range("B40").Value = "Celestia"
$ws.Range("C40").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D40").Value = "15.21"
$ws.Range("E40").Value = "  -9.91%  "

$ws.Range("E41").Value = "  -6.17%  "

$ws.Range("D42").Value = "2.25"
$ws.Range("E42").Value = "  -3.01%  "

$ws.Range("D43").Value = "1.928.61"
$ws.Range("E43").Value = "  -3.35%  "

$ws.Range("D44").Value = "0.0261"
$ws.Range("E44").Value = "  -6.80%  "

$ws.Range("D45").Value = "17.22"
$ws.Range("E45").Value = "  -7.08%  "

$ws.Range("D46").Value = "9.28"
$ws.Range("E46").Value = "  -2.15%  "

$ws.Range("E47").Value = "  -10.38%  "

$ws.Range("D48").Value = "2.538.79"
$ws.Range("E48").Value = "  -6.54%  "

$ws.Range("D49").Value = "91.63"
$ws.Range("E49").Value = "  -5.39%  "

$ws.Range("D50").Value = "69.19"
$ws.Range("E50").Value = "  -8.41%  "

$ws.Range("D51").Value = "62.56"
$ws.Range("E51").Value = "  -6.05%  "

# Reset the Price column style back to its original (no explicit number format)
# while keeping the values stored as text.
$priceRange.Style = "Normal"